$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Formatting: bold font, thin border on all sides, centered/top aligned.
# Build the style once on B1, then copy the formatting over to A2 so both
# cells land on the very same style record instead of each independently
# accumulating intermediate (unused) style entries.
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.Borders.LineStyle = 1
$r1.Borders.Weight = 2
$r1.HorizontalAlignment = -4108
$r1.VerticalAlignment = -4160

$r1.Copy()
$r2 = $ws.Range("A2")
$r2.PasteSpecial(-4122)
$excel.CutCopyMode = $false
